# Actualización automática de tasas-transfi.xlsx
# - Refresh the "Conversión del día" summary text on Hoja1!A1 with the
#   latest Binance conversion rates/amounts.
# - Refresh the raw quote cells on the "tasas" sheet that feed the
#   Binance formulas (N10/O10 = bs->pesos quote, N12/O12 = pesos->bs quote).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("tasas")

$conversionText = @"
Conversión del día 💰
✅ Dólar paralelo: 68

Binance
✅ 1000 Bs = 10.19 = 41882.58 pesos
✅ 41882.58 pesos = 10.18 = 964.88 Bs

Promedio competencia
✅ Tasa pesos: 20
✅ Tasa Bs: 20
✅ % Ganancia: 20%
"@

$ws1.Range("A1").Value = $conversionText

$ws2.Range("N10").Value = 98.179
$ws2.Range("O10").Value = 4111.99
$ws2.Range("N12").Value = 4115
$ws2.Range("O12").Value = 94.8
